$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 75, shifting existing rows 75-134 down to 76-135.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with a new weekly price record
# (same market/category metadata as the surrounding rows, new date + prices).
$ws.Range("A75").Value = 9
$ws.Range("B75").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C75").Value = "Metropolitana"
$ws.Range("D75").Value = 45096
$ws.Range("E75").Value = 13
$ws.Range("F75").Value = 100114007
$ws.Range("G75").Value = "Jengibre"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 430
$ws.Range("K75").Value = 17000
$ws.Range("L75").Value = 18000
$ws.Range("M75").Value = 17500
$ws.Range("N75").Value = "$/caja 13 kilos"
$ws.Range("O75").Value = "Perú"
$ws.Range("P75").Value = 1346
$ws.Range("Q75").Value = 13
$ws.Range("R75").Value = "Hortaliza"
